$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.008.88"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "2.051.00"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.58"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("E6").Value = "  +1.43%  "

$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.71"
$ws.Range("E7").Value = "  +2.88%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0781"
$ws.Range("E10").Value = "  +3.33%  "

$ws.Range("E11").Value = "  +1.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.81"
$ws.Range("E12").Value = "  +5.34%  "

$ws.Range("D13").Value = "2.347.40"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.61"
$ws.Range("E14").Value = "  +7.56%  "

$ws.Range("E15").Value = "  -2.95%  "

$ws.Range("D16").Value = "2.047.95"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "37.018.93"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.41"
$ws.Range("E18").Value = "  +15.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.41"
$ws.Range("E19").Value = "  +3.02%  "

$ws.Range("D20").Value = "0.0₃0896"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("E21").Value = "  +1.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.27"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  -2.10%  "

$ws.Range("E25").Value = "  +11.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.98"
$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.10"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.79"
$ws.Range("E28").Value = "  -1.43%  "

$ws.Range("E29").Value = "  +1.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.12"
$ws.Range("E30").Value = "  +6.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.69"
$ws.Range("E31").Value = "  +2.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0614"
$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.44"
$ws.Range("E33").Value = "  +2.35%  "

$ws.Range("E34").Value = "  +1.18%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.22"
$ws.Range("E36").Value = "  -2.29%  "

$ws.Range("E37").Value = "  -2.67%  "

$ws.Range("B38").Value = "Cronos"
$ws.Range("C38").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.107"
$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.34"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.12"
$ws.Range("E40").Value = "  +12.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.97"
$ws.Range("E41").Value = "  +25.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0221"
$ws.Range("E42").Value = "  -1.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.33"
$ws.Range("E43").Value = "  -5.36%  "

$ws.Range("E44").Value = "  -1.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.84"
$ws.Range("E45").Value = "  -1.19%  "

$ws.Range("E46").Value = "  +4.43%  "

$ws.Range("D47").Value = "1.282.89"
$ws.Range("E47").Value = "  -1.33%  "

$ws.Range("E48").Value = "  -2.03%  "

$ws.Range("D49").Value = "2.238.38"
$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("E50").Value = "  -1.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.49"
$ws.Range("E51").Value = "  -20.52%  "
